$wb = $excel.ActiveWorkbook

# --- Update the "Logs" sheet: append a new row of mail-log data ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(37, 1).Value = "Wat zijn jullie voorwaarden?"
$logs.Cells.Item(37, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(37, 3).Value = "Testmail #5: Wat zijn jullie voorwaarden?"
$logs.Cells.Item(37, 4).Value = "Productinformatie"
$logs.Cells.Item(37, 5).Value = "Beste klant,`nDank voor uw interesse in onze diensten. Voor informatie over onze voorwaarden kunt u terecht op onze website onder de sectie 'Algemene Voorwaarden'. Mocht u specifieke vragen hebben, dan helpen wij u graag verder. Aarzel niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item(37, 6).Value = "2025-06-26 23:19:51"
$logs.Cells.Item(37, 7).Value = "Ja"
$logs.Cells.Item(37, 8).Value = "Nee"
$logs.Cells.Item(37, 9).Value = "Ja"

# Avoid an artificial custom row-height being stamped on the new row
# (the multi-line "Antwoord" text would otherwise trigger an auto custom
# height); AutoFit brings it back to the sheet default, with no explicit
# ht/customHeight markup - matching the other rows.
$logs.Rows.Item(37).AutoFit()

# --- Extend conditional formatting ranges (D/G/H/I) to include row 37 ---
$cfRanges = @(
  @{old="D2:D36"; new="D2:D37"},
  @{old="G2:G36"; new="G2:G37"},
  @{old="H2:H36"; new="H2:H37"},
  @{old="I2:I36"; new="I2:I37"}
)

foreach ($cfr in $cfRanges) {
  $fcs = $logs.Range($cfr.old).FormatConditions
  $newRange = $logs.Range($cfr.new)
  for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($newRange)
  }
}

# --- Update the "Dashboard" sheet summary table ---
# Row 4 becomes Productinformatie (count 4), row 5 becomes Retour / Terugbetaling (count 3)
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(4, 1).Value = "Productinformatie"
$dash.Cells.Item(4, 2).Value = 4
$dash.Cells.Item(5, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(5, 2).Value = 3
